$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B213').Value = 7543661
$ws.Range('F213').Value = 'Stade Plabennecois'
$ws.Range('G213').Value = 'Locmine SaintColomban'
$ws.Range('H213').Value = 2
$ws.Range('I213').Value = 1
$ws.Range('J213').Value = 'H'
$ws.Range('K213').Value = 2.4
$ws.Range('L213').Value = 3
$ws.Range('M213').Value = 2.75
$ws.Range('N213').Value = 2.4
$ws.Range('O213').Value = 3
$ws.Range('P213').Value = 2.75
$ws.Range('Q213').Value = 0
$ws.Range('R213').Value = 1.775
$ws.Range('S213').Value = 2.025
$ws.Range('T213').Value = 2.75
$ws.Range('U213').Value = 1.925
$ws.Range('V213').Value = 1.875
$ws.Range('W213').Value = 1.4
$ws.Range('X213').Value = -1
$ws.Range('Y213').Value = -1
$ws.Range('Z213').Value = 0.7749999999999999
$ws.Range('AA213').Value = -1
$ws.Range('AB213').Value = 0.4625
$ws.Range('AC213').Value = -0.5
$ws.Range('B214').Value = 7543688
$ws.Range('F214').Value = 'Saint Priest'
$ws.Range('G214').Value = 'Vaulx en Velin'
$ws.Range('H214').Value = 5
$ws.Range('I214').Value = 0
$ws.Range('J214').Value = 'H'
$ws.Range('K214').Value = 1.571
$ws.Range('L214').Value = 3.4
$ws.Range('M214').Value = 5.5
$ws.Range('N214').Value = 1.5
$ws.Range('O214').Value = 3.6
$ws.Range('P214').Value = 5.5
$ws.Range('Q214').Value = -1
$ws.Range('R214').Value = 1.85
$ws.Range('S214').Value = 1.95
$ws.Range('T214').Value = 2.75
$ws.Range('U214').Value = 1.825
$ws.Range('V214').Value = 1.975
$ws.Range('W214').Value = 0.5
$ws.Range('Z214').Value = 0.8500000000000001
$ws.Range('AA214').Value = -1
$ws.Range('AB214').Value = 0.825
$ws.Range('AC214').Value = -1
$ws.Range('B215').Value = 7543690
$ws.Range('F215').Value = 'Feurs US'
$ws.Range('G215').Value = 'Chambery'
$ws.Range('H215').Value = 1
$ws.Range('I215').Value = 1
$ws.Range('J215').Value = 'D'
$ws.Range('K215').Value = 2.4
$ws.Range('L215').Value = 2.75
$ws.Range('M215').Value = 3
$ws.Range('N215').Value = 2.4
$ws.Range('O215').Value = 2.8
$ws.Range('P215').Value = 2.9
$ws.Range('Q215').Value = -0.25
$ws.Range('R215').Value = 2
$ws.Range('S215').Value = 1.7
$ws.Range('T215').Value = 2.25
$ws.Range('U215').Value = 1.9
$ws.Range('V215').Value = 1.9
$ws.Range('W215').Value = -1
$ws.Range('X215').Value = 1.8
$ws.Range('Y215').Value = -1
$ws.Range('Z215').Value = -0.5
$ws.Range('AA215').Value = 0.35
$ws.Range('AB215').Value = -0.5
$ws.Range('AC215').Value = 0.45
$ws.Range('B221').Value = 7547183
$ws.Range('F221').Value = 'Troyes II'
$ws.Range('G221').Value = 'Reims SteAnne'
$ws.Range('H221').Value = 2
$ws.Range('I221').Value = 1
$ws.Range('J221').Value = 'H'
$ws.Range('K221').Value = 2.5
$ws.Range('L221').Value = 2.8
$ws.Range('M221').Value = 2.8
$ws.Range('N221').Value = 2.15
$ws.Range('O221').Value = 2.875
$ws.Range('P221').Value = 3.3
$ws.Range('Q221').Value = -0.25
$ws.Range('R221').Value = 1.9
$ws.Range('S221').Value = 1.9
$ws.Range('W221').Value = 1.15
$ws.Range('X221').Value = -1
$ws.Range('Y221').Value = -1
$ws.Range('Z221').Value = 0.8999999999999999
$ws.Range('AA221').Value = -1
$ws.Range('AB221').Value = 0.4125
$ws.Range('AC221').Value = -0.5
$ws.Range('B222').Value = 7547182
$ws.Range('F222').Value = 'St Etienne II'
$ws.Range('G222').Value = 'Valence'
$ws.Range('H222').Value = 4
$ws.Range('I222').Value = 0
$ws.Range('J222').Value = 'H'
$ws.Range('K222').Value = 2.4
$ws.Range('L222').Value = 2.8
$ws.Range('M222').Value = 2.9
$ws.Range('N222').Value = 2.05
$ws.Range('O222').Value = 3
$ws.Range('R222').Value = 1.85
$ws.Range('S222').Value = 1.95
$ws.Range('T222').Value = 2.5
$ws.Range('U222').Value = 1.9
$ws.Range('V222').Value = 1.9
$ws.Range('W222').Value = 1.05
$ws.Range('Z222').Value = 0.8500000000000001
$ws.Range('AA222').Value = -1
$ws.Range('AB222').Value = 0.8999999999999999
$ws.Range('AC222').Value = -1
$ws.Range('B223').Value = 7547167
$ws.Range('F223').Value = 'Entente SCM'
$ws.Range('G223').Value = 'Marseille II'
$ws.Range('H223').Value = 1
$ws.Range('I223').Value = 1
$ws.Range('J223').Value = 'D'
$ws.Range('K223').Value = 1.909
$ws.Range('L223').Value = 3.25
$ws.Range('M223').Value = 3.5
$ws.Range('N223').Value = 1.909
$ws.Range('O223').Value = 3.25
$ws.Range('P223').Value = 3.5
$ws.Range('Q223').Value = -0.5
$ws.Range('R223').Value = 1.975
$ws.Range('S223').Value = 1.825
$ws.Range('T223').Value = 2.75
$ws.Range('U223').Value = 1.825
$ws.Range('V223').Value = 1.975
$ws.Range('W223').Value = -1
$ws.Range('X223').Value = 2.25
$ws.Range('Y223').Value = -1
$ws.Range('Z223').Value = -1
$ws.Range('AA223').Value = 0.825
$ws.Range('AB223').Value = -1
$ws.Range('AC223').Value = 0.9750000000000001
$ws.Range('B270').Value = 7718936
$ws.Range('F270').Value = 'Stade Plabennecois'
$ws.Range('G270').Value = 'Stade Pontivy'
$ws.Range('H270').Value = 5
$ws.Range('I270').Value = 3
$ws.Range('J270').Value = 'H'
$ws.Range('K270').Value = 2.6
$ws.Range('L270').Value = 3.25
$ws.Range('M270').Value = 2.375
$ws.Range('N270').Value = 1.8
$ws.Range('O270').Value = 3.6
$ws.Range('P270').Value = 3.5
$ws.Range('Q270').Value = -0.5
$ws.Range('R270').Value = 1.85
$ws.Range('S270').Value = 1.95
$ws.Range('T270').Value = 2.5
$ws.Range('U270').Value = 1.85
$ws.Range('V270').Value = 1.95
$ws.Range('W270').Value = 0.8
$ws.Range('Z270').Value = 0.8500000000000001
$ws.Range('AA270').Value = -1
$ws.Range('AB270').Value = 0.8500000000000001
$ws.Range('B271').Value = 7718935
$ws.Range('F271').Value = 'Laval II'
$ws.Range('G271').Value = 'Niort II'
$ws.Range('H271').Value = 3
$ws.Range('K271').Value = 1.571
$ws.Range('L271').Value = 3.8
$ws.Range('M271').Value = 4.75
$ws.Range('N271').Value = 1.571
$ws.Range('O271').Value = 3.8
$ws.Range('P271').Value = 4.75
$ws.Range('Q271').Value = -1
$ws.Range('R271').Value = 1.95
$ws.Range('S271').Value = 1.85
$ws.Range('T271').Value = 3
$ws.Range('W271').Value = 0.571
$ws.Range('Z271').Value = 0.95
$ws.Range('B272').Value = 7718964
$ws.Range('F272').Value = 'ASPTT Dijon'
$ws.Range('G272').Value = 'Vesoul'
$ws.Range('K272').Value = 1.833
$ws.Range('L272').Value = 3.25
$ws.Range('M272').Value = 3.8
$ws.Range('N272').Value = 1.833
$ws.Range('O272').Value = 3.25
$ws.Range('P272').Value = 3.8
$ws.Range('Q272').Value = -0.5
$ws.Range('R272').Value = 1.875
$ws.Range('S272').Value = 1.925
$ws.Range('T272').Value = 2.25
$ws.Range('U272').Value = 1.925
$ws.Range('V272').Value = 1.875
$ws.Range('W272').Value = 0.833
$ws.Range('Z272').Value = 0.875
$ws.Range('AA272').Value = -1
$ws.Range('AB272').Value = 0.925
$ws.Range('B273').Value = 7718933
$ws.Range('F273').Value = 'Lannion FC'
$ws.Range('G273').Value = 'Vannes OC'
$ws.Range('H273').Value = 2
$ws.Range('I273').Value = 1
$ws.Range('J273').Value = 'H'
$ws.Range('K273').Value = 2.25
$ws.Range('L273').Value = 3.25
$ws.Range('M273').Value = 2.75
$ws.Range('N273').Value = 2.25
$ws.Range('O273').Value = 3.25
$ws.Range('P273').Value = 2.75
$ws.Range('Q273').Value = -0.25
$ws.Range('R273').Value = 2.025
$ws.Range('S273').Value = 1.775
$ws.Range('T273').Value = 2.5
$ws.Range('U273').Value = 1.9
$ws.Range('V273').Value = 1.9
$ws.Range('W273').Value = 1.25
$ws.Range('Z273').Value = 0.95
$ws.Range('AA273').Value = -1
$ws.Range('AB273').Value = 0.8999999999999999
$ws.Range('G320').Value = 'Vannes OC'
$ws.Range('B334').Value = 7874249
$ws.Range('F334').Value = 'Blagnac'
$ws.Range('G334').Value = 'Anglet Genets'
$ws.Range('H334').Value = 0
$ws.Range('I334').Value = 0
$ws.Range('J334').Value = 'D'
$ws.Range('K334').Value = 2.1
$ws.Range('L334').Value = 3.2
$ws.Range('M334').Value = 3.1
$ws.Range('N334').Value = 2.1
$ws.Range('O334').Value = 3.2
$ws.Range('P334').Value = 3.1
$ws.Range('Q334').Value = -0.25
$ws.Range('R334').Value = 1.875
$ws.Range('S334').Value = 1.925
$ws.Range('T334').Value = 2.25
$ws.Range('U334').Value = 1.925
$ws.Range('V334').Value = 1.875
$ws.Range('W334').Value = -1
$ws.Range('X334').Value = 2.2
$ws.Range('Y334').Value = -1
$ws.Range('Z334').Value = -0.5
$ws.Range('AA334').Value = 0.4625
$ws.Range('AB334').Value = -1
$ws.Range('AC334').Value = 0.875
$ws.Range('B335').Value = 7874234
$ws.Range('F335').Value = 'Montpellier II'
$ws.Range('G335').Value = 'Marseille II'
$ws.Range('K335').Value = 1.8
$ws.Range('L335').Value = 3.2
$ws.Range('M335').Value = 4
$ws.Range('N335').Value = 1.8
$ws.Range('O335').Value = 3.2
$ws.Range('P335').Value = 4
$ws.Range('Q335').Value = -0.5
$ws.Range('R335').Value = 1.85
$ws.Range('S335').Value = 1.95
$ws.Range('T335').Value = 2.75
$ws.Range('U335').Value = 1.95
$ws.Range('V335').Value = 1.85
$ws.Range('Z335').Value = -1
$ws.Range('AA335').Value = 0.95
$ws.Range('AB335').Value = -1
$ws.Range('AC335').Value = 0.8500000000000001
$ws.Range('B336').Value = 7874397
$ws.Range('F336').Value = 'Lille II'
$ws.Range('G336').Value = 'Compiegne'
$ws.Range('H336').Value = 3
$ws.Range('I336').Value = 1
$ws.Range('J336').Value = 'H'
$ws.Range('K336').Value = 1.5
$ws.Range('L336').Value = 4
$ws.Range('M336').Value = 5
$ws.Range('N336').Value = 1.5
$ws.Range('O336').Value = 4
$ws.Range('P336').Value = 4.75
$ws.Range('Q336').Value = -1
$ws.Range('R336').Value = 1.875
$ws.Range('S336').Value = 1.925
$ws.Range('W336').Value = 0.5
$ws.Range('X336').Value = -1
$ws.Range('Y336').Value = -1
$ws.Range('Z336').Value = 0.875
$ws.Range('AA336').Value = -1
$ws.Range('AB336').Value = 0.95
$ws.Range('AC336').Value = -1
$ws.Range('B337').Value = 7874399
$ws.Range('F337').Value = 'Dijon II'
$ws.Range('G337').Value = 'Selongey'
$ws.Range('H337').Value = 1
$ws.Range('I337').Value = 1
$ws.Range('J337').Value = 'D'
$ws.Range('K337').Value = 1.4
$ws.Range('L337').Value = 4
$ws.Range('M337').Value = 6.5
$ws.Range('N337').Value = 1.571
$ws.Range('O337').Value = 3.75
$ws.Range('R337').Value = 1.95
$ws.Range('S337').Value = 1.85
$ws.Range('T337').Value = 2.75
$ws.Range('U337').Value = 1.825
$ws.Range('V337').Value = 1.975
$ws.Range('W337').Value = -1
$ws.Range('X337').Value = 2.75
$ws.Range('Y337').Value = -1
$ws.Range('Z337').Value = -1
$ws.Range('AA337').Value = 0.8500000000000001
$ws.Range('AB337').Value = -1
$ws.Range('AC337').Value = 0.9750000000000001
